$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying analysis script was changed to omit "Other" land from the
# statistics. This removes the "Otherland_pcnt" row (originally row 5) from
# the summary table entirely; all rows below it shift up by one, and the
# weighted percentage / count statistics are recomputed without "Other"
# land contributing to the totals.

# Delete the "Otherland_pcnt" row (row 5). Excel will automatically drop the
# now-unused "Otherland_pcnt" entry from the shared string table and shift
# the remaining rows up.
$ws.Rows.Item(5).Delete()

# Recomputed weighted percentages / counts for the affected rows, now that
# "Other" land is excluded from the totals.

# Row 2: Cropland_pcnt
$ws.Range("B2").Value = 25.679052352905273
$ws.Range("C2").Value = 24.366315841674805
$ws.Range("D2").Value = 22.934646606445313
$ws.Range("E2").Value = 22.422433853149414
$ws.Range("F2").Value = 21.424379348754883
$ws.Range("G2").Value = 21.098880767822266
$ws.Range("H2").Value = 21.753026962280273

# Row 3: CRPland_pcnt
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0.96704113483428955
$ws.Range("D3").Value = 2.3926632404327393
$ws.Range("E3").Value = 2.3010172843933105
$ws.Range("F3").Value = 2.2188904285430908
$ws.Range("G3").Value = 2.3002748489379883
$ws.Range("H3").Value = 1.6927685737609863

# Row 4: Forestland_pcnt
$ws.Range("B4").Value = 28.656833648681641
$ws.Range("C4").Value = 28.828330993652344
$ws.Range("D4").Value = 28.879106521606445
$ws.Range("E4").Value = 28.993467330932617
$ws.Range("F4").Value = 29.069923400878906
$ws.Range("G4").Value = 29.101343154907227
$ws.Range("H4").Value = 29.170219421386719

# Row 5: Pastureland_pcnt (was old row 6)
$ws.Range("B5").Value = 12.877095222473145
$ws.Range("C5").Value = 12.968668937683105
$ws.Range("D5").Value = 12.631027221679688
$ws.Range("E5").Value = 12.455606460571289
$ws.Range("F5").Value = 12.846061706542969
$ws.Range("G5").Value = 12.666543960571289
$ws.Range("H5").Value = 12.418011665344238

# Row 6: Rangeland_pcnt (was old row 7)
$ws.Range("B6").Value = 29.254037857055664
$ws.Range("C6").Value = 28.917793273925781
$ws.Range("D6").Value = 28.710187911987305
$ws.Range("E6").Value = 28.638595581054688
$ws.Range("F6").Value = 28.64451789855957
$ws.Range("G6").Value = 28.635763168334961
$ws.Range("H6").Value = 28.558055877685547

# Row 7: Urbanland_pcnt (was old row 8)
$ws.Range("B7").Value = 3.5329799652099609
$ws.Range("C7").Value = 3.9518492221832275
$ws.Range("D7").Value = 4.4523677825927734
$ws.Range("E7").Value = 5.188880443572998
$ws.Range("F7").Value = 5.7962260246276855
$ws.Range("G7").Value = 6.1971955299377441
$ws.Range("H7").Value = 6.4079184532165527

# Rows 8-19 (lccL1_pcnt ... lccL78_pcnt, were old rows 9-20) are unaffected
# by the removal of "Other" land and keep their original values.

# Row 20: crop_nr (was old row 21) - recomputed counts
$ws.Range("B20").Value = 58.686473846435547
$ws.Range("C20").Value = 105.58715057373047
$ws.Range("D20").Value = 124.09786224365234
$ws.Range("E20").Value = 119.66957855224609
$ws.Range("F20").Value = 84.831069946289063
$ws.Range("G20").Value = 89.905418395996094
$ws.Range("H20").Value = 118.20491027832031

# Row 21: forest_nr (was old row 22) - recomputed counts; B21:E21 stay blank
$ws.Range("F21").Value = 19.238698959350586
$ws.Range("G21").Value = 17.943864822387695
$ws.Range("H21").Value = 15.511234283447266

# Row 22: urban_nr (was old row 23) - recomputed counts; B22:D22 stay blank
$ws.Range("E22").Value = 26516.9375
$ws.Range("F22").Value = 33384.97265625
$ws.Range("G22").Value = 53524.9765625
$ws.Range("H22").Value = 50263.3671875
